$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the "Work out openness clse ness..." paragraph and blank
#    its text down to a single space, dropping the old _GoBack bookmark
#    that used to live there.
# ------------------------------------------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Work out openness clse ness of doors*") {
        $target1 = $p
        break
    }
}

if ($target1 -ne $null) {
    $r1 = $target1.Range
    $txtRange = $d.Range($r1.Start, $r1.End - 1)
    $txtRange.Text = " "
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "I notice in ravenrock..." paragraph, insert a new
#    list paragraph right after it with the new bug note, and restore
#    the _GoBack bookmark at the end of that new paragraph's text
#    (right after the run, before the paragraph mark).
# ------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I notice in ravenrock*") {
        $target2 = $p
        break
    }
}

if ($target2 -ne $null) {
    $target2.Range.InsertParagraphAfter()

    $newPara = $target2.Next()
    $apos = [char]0x2019
    $newRange = $newPara.Range
    $newRange.Text = "If you click an opening door, it loses it" + $apos + "s mind"

    # A zero-width Range sitting exactly at "end of paragraph text" can
    # resolve incorrectly when handed straight to Bookmarks.Add, so we
    # park a throwaway character there first, anchor the bookmark just
    # before it (now a safe, non-boundary position), then remove the
    # throwaway character again.
    $insertPos = $newPara.Range.End - 1
    $insRange = $d.Range($insertPos, $insertPos)
    $insRange.InsertAfter("X")

    $bmRange = $d.Range($insertPos, $insertPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $delRange = $d.Range($insertPos, $insertPos + 1)
    $delRange.Delete()
}
